$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "11/14/2025"
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").Value = 507.4490000000005
$ws.Range("C14").Value = 0.04877337427012365
$ws.Range("D14").Value = 25
